# Q3 Update - 2025
# Applies the quarterly data refresh to the UN-SUR sheet:
#  - refreshes the short-url token used throughout the sheet
#  - updates a handful of asylum_seekers (O) figures for 2024 rows
#  - drops the last two detail rows (Stateless / Suriname->Suriname) and
#    folds the remaining Venezuela->Suriname row up into row 55 with
#    refreshed figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (old rows 55 "Stateless" and 56
# "Suriname -> Suriname"); this shifts the old row 57
# (Venezuela -> Suriname) up to become the new row 55.
$ws.Rows("55:56").Delete()

# The "items" sequence id on the folded-up row keeps its original value.
$ws.Range("D55").Value = "54"

# Refreshed asylum_seekers (column O) figures.
$ws.Range("O55").Value = "482"
$ws.Range("T55").Value = "0"

$ws.Range("O48").Value = "5"
$ws.Range("O49").Value = "5"
$ws.Range("O53").Value = "5"
$ws.Range("O54").Value = "5"

$ws.Range("O50").Value = "2705"

# Refresh the short-url token shared by every data row.
$ws.Range("B2:B55").Value = "L82Ktu"
